# TradingPL.xlsx - "made chagnes for month append"
# Append a new "May-24" month block (Qty / Value / Rate) in columns E:G,
# mirroring the existing "Apr-24" block that lives in columns B:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column widths: the existing sheet defines widths one column past the
# last column that actually holds data (col E already had a width even
# though only A:D had data). Keep that same convention for the new
# block: F, G, H get the standard data-column width.
# ---------------------------------------------------------------------
$ws.Range("F1").EntireColumn.ColumnWidth = 15
$ws.Range("G1").EntireColumn.ColumnWidth = 15
$ws.Range("H1").EntireColumn.ColumnWidth = 15

# ---------------------------------------------------------------------
# Row 1: month header ("May-24") plus two blank cells to the right of it
# (mirrors B1="Apr-24", C1="", D1="").
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "May-24"

# ---------------------------------------------------------------------
# Row 2: column headers for the new month block (mirrors B2:D2).
# ---------------------------------------------------------------------
$ws.Range("E2").Value = "Qty"
$ws.Range("F2").Value = "Value"
$ws.Range("G2").Value = "Rate"

# ---------------------------------------------------------------------
# Helper data: rows that carry real figures. Qty/Value are numbers,
# Rate is stored as text (same convention as the existing D column).
# ---------------------------------------------------------------------

# Row 5: Sales-Mono Shade Net
$ws.Range("E5").Value = 43
$ws.Range("F5").Value = 10400
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "241.86"

# Row 7: Sales-Tape Shade Net
$ws.Range("E7").Value = 298.25
$ws.Range("F7").Value = 70710
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "237.08"

# Row 8: Sales-Weed Mate Fabrics
$ws.Range("E8").Value = 949.42
$ws.Range("F8").Value = 183407.9
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "193.18"

# Row 9: Sales-PP Woven Sacks
$ws.Range("E9").Value = 6605
$ws.Range("F9").Value = 978441
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "148.14"

# Row 12: Add: Purchase MSN (Qty only)
$ws.Range("E12").Value = 549

# Row 13: Add: Purchase PP Sacks
$ws.Range("E13").Value = 5918
$ws.Range("F13").Value = 759300
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "128.30"

# Row 14: Add: Purchase TSN
$ws.Range("E14").Value = 298.25
$ws.Range("F14").Value = 56048
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "187.92"

# ---------------------------------------------------------------------
# All remaining rows in the new E:G block are blank, like their B:D
# counterparts. Mark them present (and keep them text-typed, matching
# the Rate column convention) without disturbing the numeric cells
# above.
# ---------------------------------------------------------------------
$blankRows = @(3,4,6,10,11,15,16,17,18,19,20,21,22,23)
foreach ($r in $blankRows) {
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("G$r").NumberFormat = "@"
}
# Row 12 only needs F and G blanked (E12 already holds the Qty value).
$ws.Range("F12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"
